# Update the "Velocity" workbook:
#  - B4 changes from "=(20+4)/2" (12) to "=1+3+5" (9) -- the velocity total
#  - the worksheet's print orientation is explicitly set to portrait
#  - the user's selection ends up on B5 (the cell just below the edited one)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated velocity": recompute the total with a new formula
$ws.Range("B4").Formula = "=1+3+5"

# "added the project summary section to the report": turn on an explicit
# (portrait) page setup for printing the report
$ws.PageSetup.Orientation = 1

# leave the cursor on the cell below the one that was just edited
$ws.Range("B5").Select() | Out-Null
